# Convert the Week_Start / Week_End date columns into a single "Week"
# text column formatted as "m/d/yyyy - m/d/yyyy", and drop the old
# Week_End column so Total_Positive shifts left into column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 17

# Step 1: build the combined "8/7/2020 - 8/13/2020" style strings in a
# scratch column (D) from the existing Week_Start (A) / Week_End (B) dates.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Range("D$r").Formula = "=TEXT(A$r,""m/d/yyyy"")&"" - ""&TEXT(B$r,""m/d/yyyy"")"
}

# Step 2: copy those computed strings back onto column A as literal values.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Range("A$r").Value = $ws.Range("D$r").Value2
}

# scratch column no longer needed
$ws.Range("D$firstRow`:D$lastRow").Clear()

# Step 3: drop the old Week_End column (B) - Total_Positive (old C) slides
# left into column B.
$ws.Columns.Item(2).Delete()

# Step 4: headers
$ws.Range("A1").Value = "Week"
$ws.Range("B1").Value = "Total_Positive"

# Step 5: give the new Week text column a date-flavoured number format
# (matches the workbook author's original column formatting) while
# keeping a single shared style record for every data row.
$ws.Range("A$firstRow").NumberFormat = "mm-dd-yy"
$ws.Range("A$firstRow").Copy()
$ws.Range("A$($firstRow+1):A$lastRow").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Step 6: the Week column now holds much longer text, so widen it to fit;
# the old fixed-width Total_Positive column reverts to the sheet default.
$ws.Columns.Item(1).AutoFit()

# Step 7: reflect the author's view state (zoomed in, selection resting
# on the data that was just retyped).
$excel.ActiveWindow.Zoom = 168
$ws.Range("C5").Select()
